# Added the circuitry for the ALU
# Adds 5 new columns (U:isJump, V:jmpInterest, W:jmpCond, X:clrFlags, Y:setFlag)
# to the microinstruction table and fills in placeholder content to match the
# existing table layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) -----------------------------------------------
# Write in this exact order so new shared strings land in the expected
# slots (isJump, jmpInterest, jmpCond, 0000, clrFlags, setFlag, 1111).
$ws.Range("U1").Value = "isJump"
$ws.Range("V1").Value = "jmpInterest"
$ws.Range("W1").Value = "jmpCond"

# First "0000" placeholder write (seeds the shared string before clrFlags/setFlag)
$ws.Range("V3").NumberFormat = "@"
$ws.Range("V3").Value = "0000"

$ws.Range("X1").Value = "clrFlags"
$ws.Range("Y1").Value = "setFlag"

# --- Column width for the new V column ---------------------------------
$ws.Range("V1").ColumnWidth = 10.75

# --- Row 2: single styled placeholder cell -----------------------------
$ws.Range("W2").NumberFormat = "@"

# --- Data rows: U/X numeric 0, V/W/Y text "0000" (Y22 is the exception) -
$dataRows = @(3,5,7,9,11,13,15,17,19,21,22,23)
foreach ($r in $dataRows) {
    $ws.Cells.Item($r, 21).Value = 0          # U
    $ws.Cells.Item($r, 22).NumberFormat = "@" # V
    $ws.Cells.Item($r, 22).Value = "0000"
    $ws.Cells.Item($r, 23).NumberFormat = "@" # W
    $ws.Cells.Item($r, 23).Value = "0000"
    $ws.Cells.Item($r, 24).Value = 0          # X
    $ws.Cells.Item($r, 25).NumberFormat = "@" # Y
    if ($r -eq 22) {
        $ws.Cells.Item($r, 25).Value = "1111"
    } else {
        $ws.Cells.Item($r, 25).Value = "0000"
    }
}

# --- Header/group rows: V/W/Y styled but empty --------------------------
$headerRows = @(4,6,8,10,12,14,16,18,20)
foreach ($r in $headerRows) {
    $ws.Cells.Item($r, 22).NumberFormat = "@" # V
    $ws.Cells.Item($r, 23).NumberFormat = "@" # W
    $ws.Cells.Item($r, 25).NumberFormat = "@" # Y
}

# --- Filler rows 24-35: V/W/Y styled but empty ---------------------------
for ($r = 24; $r -le 35; $r++) {
    $ws.Cells.Item($r, 22).NumberFormat = "@" # V
    $ws.Cells.Item($r, 23).NumberFormat = "@" # W
    $ws.Cells.Item($r, 25).NumberFormat = "@" # Y
}

# --- Row 36: W/Y styled but empty (no V) ---------------------------------
$ws.Cells.Item(36, 23).NumberFormat = "@" # W
$ws.Cells.Item(36, 25).NumberFormat = "@" # Y

# --- Row 37 (new row): W/Y styled but empty ------------------------------
$ws.Cells.Item(37, 23).NumberFormat = "@" # W
$ws.Cells.Item(37, 25).NumberFormat = "@" # Y

# --- Rows 38-40 (new rows): Y styled but empty ---------------------------
for ($r = 38; $r -le 40; $r++) {
    $ws.Cells.Item($r, 25).NumberFormat = "@" # Y
}

# --- Selection matches the saved view in the target workbook ------------
$ws.Range("S22").Select() | Out-Null
